$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $text)
    $range.Formula = '="' + $text + '"'
    $range.Copy()
    $range.PasteSpecial(-4163)
}

Set-TextValue $ws.Range("D2") "43.020.47"
$ws.Range("E2").Value = "  -5.01%  "

Set-TextValue $ws.Range("D3") "2.224.96"
$ws.Range("E3").Value = "  -6.05%  "

$ws.Range("E4").Value = "  +0.05%  "

Set-TextValue $ws.Range("D5") "318.36"
$ws.Range("E5").Value = "  +2.45%  "

$ws.Range("E6").Value = "  -7.45%  "

Set-TextValue $ws.Range("D7") "0.592"
$ws.Range("E7").Value = "  -6.14%  "

$ws.Range("E8").Value = "  +0.01%  "

Set-TextValue $ws.Range("D9") "0.564"
$ws.Range("E9").Value = "  -7.01%  "

Set-TextValue $ws.Range("D10") "37.42"
$ws.Range("E10").Value = "  -8.42%  "

Set-TextValue $ws.Range("D11") "54.12"
$ws.Range("E11").Value = "  -2.68%  "

Set-TextValue $ws.Range("D12") "0.0834"

Set-TextValue $ws.Range("D13") "7.83"
$ws.Range("E13").Value = "  -7.30%  "

$ws.Range("E14").Value = "  -2.44%  "

Set-TextValue $ws.Range("D15") "0.867"
$ws.Range("E15").Value = "  -11.14%  "

Set-TextValue $ws.Range("D16") "2.562.27"
$ws.Range("E16").Value = "  -6.04%  "

Set-TextValue $ws.Range("D17") "14.31"
$ws.Range("E17").Value = "  -5.85%  "

Set-TextValue $ws.Range("D18") "2.219.53"
$ws.Range("E18").Value = "  -6.22%  "

Set-TextValue $ws.Range("D19") "42.922.09"
$ws.Range("E19").Value = "  -5.05%  "

Set-TextValue $ws.Range("D20") "15.02"
$ws.Range("E20").Value = "  +3.59%  "

Set-TextValue $ws.Range("D21") "0.0₃0966"
$ws.Range("E21").Value = "  -8.71%  "

Set-TextValue $ws.Range("D22") "6.44"
$ws.Range("E22").Value = "  -10.67%  "

Set-TextValue $ws.Range("D23") "65.68"
$ws.Range("E23").Value = "  -10.10%  "

$ws.Range("E24").Value = "  -9.22%  "

Set-TextValue $ws.Range("D25") "236.78"
$ws.Range("E25").Value = "  -8.73%  "

$ws.Range("E26").Value = "  -6.64%  "

$ws.Range("E27").Value = "  +0.23%  "

Set-TextValue $ws.Range("D28") "10.07"
$ws.Range("E28").Value = "  -9.14%  "

$ws.Range("E29").Value = "  -4.63%  "

Set-TextValue $ws.Range("D30") "6.37"
$ws.Range("E30").Value = "  -11.42%  "

Set-TextValue $ws.Range("D31") "0.0902"
$ws.Range("E31").Value = "  -6.82%  "

Set-TextValue $ws.Range("D32") "20.60"
$ws.Range("E32").Value = "  -7.82%  "

Set-TextValue $ws.Range("D33") "34.49"
$ws.Range("E33").Value = "  -7.43%  "

Set-TextValue $ws.Range("D34") "157.08"
$ws.Range("E34").Value = "  -6.48%  "

$ws.Range("E35").Value = "  -6.85%  "

Set-TextValue $ws.Range("D36") "3.18"
$ws.Range("E36").Value = "  +9.98%  "

$ws.Range("E37").Value = "  +11.77%  "

Set-TextValue $ws.Range("D38") "0.122"
$ws.Range("E38").Value = "  -6.12%  "

$ws.Range("E39").Value = "  -3.60%  "

Set-TextValue $ws.Range("D40") "3.90"
$ws.Range("E40").Value = "  -1.52%  "

$ws.Range("E41").Value = "  -9.55%  "

$ws.Range("E42").Value = "  -6.84%  "

Set-TextValue $ws.Range("D43") "1.931.02"
$ws.Range("E43").Value = "  +1.64%  "

$ws.Range("E44").Value = "  +0.07%  "

Set-TextValue $ws.Range("D45") "12.62"
$ws.Range("E45").Value = "  -1.70%  "

Set-TextValue $ws.Range("D46") "89.35"
$ws.Range("E46").Value = "  -10.74%  "

Set-TextValue $ws.Range("D47") "0.210"
$ws.Range("E47").Value = "  -8.49%  "

$ws.Range("E48").Value = "  -3.98%  "

Set-TextValue $ws.Range("D49") "77.26"
$ws.Range("E49").Value = "  -6.68%  "

Set-TextValue $ws.Range("D50") "60.74"
$ws.Range("E50").Value = "  -12.53%  "

# Row 51: SEI -> Aave
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws.Range("D51") "103.80"
$ws.Range("E51").Value = "  -5.93%  "
